$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for existing rows 2-28
# from serial 45446 (2024-06-03) to serial 45447 (2024-06-04)
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = 45447
}

# Append a new data row (row 29) for case "A 21946-2024"
$ws.Range("A29").Value = "A 21946-2024"
$ws.Range("B29").Value = 45443
$ws.Range("C29").Value = 45447
$ws.Range("D29").Value = "OKÄNT"
$ws.Range("E29").Value = "OKÄNT"
$ws.Range("G29").Value = 1.7
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0

# R29 mirrors the other rows' empty, word-wrapped "Artnamn" cell
$ws.Range("R29").Value = ""
$ws.Range("R29").WrapText = $true

# Give the date columns their expected date format/style (matches other rows)
$ws.Range("B29:C29").NumberFormat = "YYYY-MM-DD"

# Row 28 picks up an explicit custom row height once row 29 is appended
$ws.Rows.Item(28).RowHeight = 15
